# Update the 2022 population figures (column F) with the refreshed data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "F2"  = 18157
    "F3"  = 80782
    "F4"  = 114101
    "F5"  = 112566
    "F6"  = 115697
    "F7"  = 278572
    "F8"  = 293800
    "F9"  = 281965
    "F10" = 234356
    "F11" = 174706
    "F12" = 93580
    "F13" = 43971
    "F14" = 18488
    "F15" = 83644
    "F16" = 119030
    "F17" = 117496
    "F18" = 118852
    "F19" = 278848
    "F20" = 280289
    "F21" = 261639
    "F22" = 205547
    "F23" = 140120
    "F24" = 69260
    "F25" = 24792
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

# Reflect the author's last selected cell when the workbook was saved.
$ws.Range("E7").Select()
